$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 44515
$ws.Range("B7").Value = 1711
$ws.Range("C7").Value = 0.03
$ws.Range("D7").Value = 187

$ws.Range("A8").Value = 44515
$ws.Range("B8").Value = 3033
$ws.Range("C8").Value = -0.04
$ws.Range("D8").Value = -243

$ws.Range("A9").Value = 44515
$ws.Range("B9").Value = 3189
$ws.Range("C9").Value = 0.05
$ws.Range("D9").Value = 297

$ws.Range("A7:A9").NumberFormat = $ws.Range("A6").NumberFormat
